# Update gh-pages output (广州-漫展信息.xlsx) to the values generated at 456a3b4.
#
# Sheet layout: 1=展览 (exhibitions), 2=演出 (performances),
#               3=本地生活 (local life), 4=全部类型 (all types, union of 1-3)
#
# Semantic change:
#   - The "Tielle" concert listing (bilibili id=90812) was removed from the
#     source feed, so its row is deleted from 演出 (sheet 2) and from the
#     merged 全部类型 (sheet 4); subsequent rows shift up by one and the
#     sheet dimension shrinks by one row on both sheets.
#   - Several events' "想去人数" (want-to-go count, column F) were
#     incremented on refresh; update those across every sheet that lists
#     the event (展览, 本地生活, 全部类型).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsShow    = $wb.Worksheets.Item(2)   # 演出
$wsLocal   = $wb.Worksheets.Item(3)   # 本地生活
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# --- 1) Remove the cancelled "Tielle" concert row -------------------------
# 演出 (sheet 2): row 2 is the Tielle listing -> delete, rows shift up.
$wsShow.Rows.Item(2).Delete()

# 全部类型 (sheet 4): row 5 is the same Tielle listing -> delete.
$wsAll.Rows.Item(5).Delete()

# --- 2) Refresh "want-to-go" counts (column F) -----------------------------

# 展览 (sheet 1)
$wsExhibit.Range("F3").Value  = 211
$wsExhibit.Range("F4").Value  = 598
$wsExhibit.Range("F5").Value  = 236
$wsExhibit.Range("F6").Value  = 392
$wsExhibit.Range("F7").Value  = 527
$wsExhibit.Range("F8").Value  = 206
$wsExhibit.Range("F9").Value  = 57
$wsExhibit.Range("F10").Value = 328
$wsExhibit.Range("F11").Value = 116
$wsExhibit.Range("F12").Value = 551
$wsExhibit.Range("F14").Value = 1706
$wsExhibit.Range("F15").Value = 292
$wsExhibit.Range("F16").Value = 1549
$wsExhibit.Range("F17").Value = 229
$wsExhibit.Range("F18").Value = 476
$wsExhibit.Range("F19").Value = 25
$wsExhibit.Range("F21").Value = 126

# 本地生活 (sheet 3)
$wsLocal.Range("F2").Value = 5237
$wsLocal.Range("F3").Value = 292

# 全部类型 (sheet 4) - rows already shifted up by one from the deletion above
$wsAll.Range("F3").Value  = 5237
$wsAll.Range("F4").Value  = 292
$wsAll.Range("F6").Value  = 211
$wsAll.Range("F13").Value = 598
$wsAll.Range("F16").Value = 236
$wsAll.Range("F17").Value = 392
$wsAll.Range("F18").Value = 527
$wsAll.Range("F19").Value = 206
$wsAll.Range("F20").Value = 57
$wsAll.Range("F22").Value = 328
$wsAll.Range("F23").Value = 116
$wsAll.Range("F26").Value = 551
$wsAll.Range("F29").Value = 1706
$wsAll.Range("F30").Value = 292
$wsAll.Range("F31").Value = 1549
$wsAll.Range("F33").Value = 229
$wsAll.Range("F34").Value = 476
$wsAll.Range("F35").Value = 25
$wsAll.Range("F38").Value = 126
